$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1949
$ws.Range("I70").Value = 2349
$ws.Range("J70").Value = 1749
$ws.Range("K70").Value = 7047
$ws.Range("L70").Value = 5247
$ws.Range("M70").Value = -6777
$ws.Range("N70").Value = -5787
$ws.Range("H73").Value = 1949
$ws.Range("I73").Value = 2349
$ws.Range("J73").Value = 1749
$ws.Range("K73").Value = 7047
$ws.Range("L73").Value = 5247
$ws.Range("M73").Value = -6111
$ws.Range("N73").Value = -7119
$ws.Range("H74").Value = 7531291.5
$ws.Range("J74").Value = 13788.471
$ws.Range("L74").Value = 13788.471
$ws.Range("N74").Value = -15660.471
$ws.Range("H77").Value = 7531291.5
$ws.Range("J77").Value = 13788.471
$ws.Range("L77").Value = 68942.355
$ws.Range("N77").Value = -78302.355
$ws.Range("H80").Value = 1437.3334
$ws.Range("I80").Value = 699.5
$ws.Range("K80").Value = 2098.5
$ws.Range("M80").Value = -1100.5
$ws.Range("H83").Value = 1437.3334
$ws.Range("I83").Value = 699.5
$ws.Range("K83").Value = 6295.5
$ws.Range("M83").Value = -1303.5
$ws.Range("H106").Value = 4499
$ws.Range("I106").Value = 4540.5835
$ws.Range("K106").Value = 4540.5835
$ws.Range("M106").Value = -3909.5835
$ws.Range("H113").Value = 52999.668
$ws.Range("I113").Value = 9000
$ws.Range("J113").Value = 74999.5
$ws.Range("K113").Value = 9000
$ws.Range("L113").Value = 74999.5
$ws.Range("M113").Value = -5746
$ws.Range("N113").Value = -81507.5
$ws.Range("H125").Value = 7791.5
$ws.Range("I125").Value = 8278
$ws.Range("K125").Value = 74502
$ws.Range("M125").Value = -72042
$ws.Range("H132").Value = 758019.25
$ws.Range("I132").Value = 803832.9399999999
$ws.Range("K132").Value = 2411498.82
$ws.Range("M132").Value = -2408968.82
$ws.Range("H136").Value = 148989.5
$ws.Range("J136").Value = 148989.5
$ws.Range("L136").Value = 148989.5
$ws.Range("N136").Value = -159189.5

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1685594
$ws.Range("I32").Value = 707.375
$ws.Range("J32").Value = 6178625
$ws.Range("K32").Value = 707.375
$ws.Range("L32").Value = 6178625
$ws.Range("M32").Value = -420.375
$ws.Range("N32").Value = -6179199
$ws.Range("H97").Value = 1682.8
$ws.Range("I97").Value = 1548.4286
$ws.Range("J97").Value = 1996.3334
$ws.Range("K97").Value = 1548.4286
$ws.Range("L97").Value = 1996.3334
$ws.Range("M97").Value = -1052.4286
$ws.Range("N97").Value = -2988.3334
$ws.Range("H122").Value = 4944.619
$ws.Range("I122").Value = 4386.643
$ws.Range("J122").Value = 6060.5713
$ws.Range("K122").Value = 13159.929
$ws.Range("L122").Value = 18181.7139
$ws.Range("M122").Value = -10709.929
$ws.Range("N122").Value = -23081.7139
$ws.Range("H132").Value = 1563731.8
$ws.Range("I132").Value = 2165706.2
$ws.Range("K132").Value = 6497118.600000001
$ws.Range("M132").Value = -6494588.600000001
$ws.Range("H135").Value = 70429
$ws.Range("J135").Value = 70429
$ws.Range("L135").Value = 70429
$ws.Range("N135").Value = -80569

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 937287.75
$ws.Range("I134").Value = 1038162.5
$ws.Range("K134").Value = 3114487.5
$ws.Range("M134").Value = -3111952.5

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13894083
$ws.Range("I99").Value = 27781388
$ws.Range("K99").Value = 27781388
$ws.Range("M99").Value = -27779890
$ws.Range("H122").Value = 2739.5454
$ws.Range("I122").Value = 2517.7144
$ws.Range("K122").Value = 7553.1432
$ws.Range("M122").Value = -5103.1432
$ws.Range("H126").Value = 13894083
$ws.Range("I126").Value = 27781388
$ws.Range("K126").Value = 83344164
$ws.Range("M126").Value = -83341694
$ws.Range("H134").Value = 37047484
$ws.Range("I134").Value = 142861710
$ws.Range("J134").Value = 12501.75
$ws.Range("K134").Value = 428585130
$ws.Range("L134").Value = 37505.25
$ws.Range("M134").Value = -428582595
$ws.Range("N134").Value = -42575.25

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 295752.34
$ws.Range("I5").Value = 1741.6666
$ws.Range("K5").Value = 5224.9998
$ws.Range("M5").Value = -5112.9998
$ws.Range("H37").Value = 219993.17
$ws.Range("J37").Value = 219993.17
$ws.Range("L37").Value = 659979.51
$ws.Range("N37").Value = -660203.51
$ws.Range("H80").Value = 1849.5
$ws.Range("I80").Value = 1799
$ws.Range("J80").Value = 1900
$ws.Range("K80").Value = 5397
$ws.Range("L80").Value = 5700
$ws.Range("M80").Value = -4461
$ws.Range("N80").Value = -7572
$ws.Range("H83").Value = 1849.5
$ws.Range("I83").Value = 1799
$ws.Range("J83").Value = 1900
$ws.Range("K83").Value = 16191
$ws.Range("L83").Value = 17100
$ws.Range("M83").Value = -11511
$ws.Range("N83").Value = -26460
$ws.Range("H107").Value = 2670.2666
$ws.Range("I107").Value = 404.5
$ws.Range("J107").Value = 4180.778
$ws.Range("K107").Value = 1213.5
$ws.Range("L107").Value = 12542.334
$ws.Range("M107").Value = 706.5
$ws.Range("N107").Value = -16382.334
$ws.Range("H131").Value = 45616004
$ws.Range("I131").Value = 66667016
$ws.Range("J131").Value = 30306172
$ws.Range("K131").Value = 200001048
$ws.Range("L131").Value = 90918516
$ws.Range("M131").Value = -199996008
$ws.Range("N131").Value = -90928596
$ws.Range("H135").Value = 295752.34
$ws.Range("I135").Value = 1741.6666
$ws.Range("K135").Value = 15674.9994
$ws.Range("M135").Value = -13139.9994

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7218.6924
$ws.Range("I102").Value = 6736.75
$ws.Range("J102").Value = 7989.8
$ws.Range("K102").Value = 6736.75
$ws.Range("L102").Value = 7989.8
$ws.Range("M102").Value = -5114.75
$ws.Range("N102").Value = -11233.8
$ws.Range("H122").Value = 7684.6665
$ws.Range("I122").Value = 5568.8335
$ws.Range("J122").Value = 11916.333
$ws.Range("K122").Value = 16706.5005
$ws.Range("L122").Value = 35748.999
$ws.Range("M122").Value = -14256.5005
$ws.Range("N122").Value = -40648.999
$ws.Range("H126").Value = 54552764
$ws.Range("I126").Value = 125002856
$ws.Range("J126").Value = 14295571
$ws.Range("K126").Value = 375008568
$ws.Range("L126").Value = 42886713
$ws.Range("M126").Value = -375006098
$ws.Range("N126").Value = -42891653
$ws.Range("H132").Value = 76924880
$ws.Range("I132").Value = 76924880
$ws.Range("K132").Value = 230774640
$ws.Range("M132").Value = -230772110

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9039.875
$ws.Range("I7").Value = 5642.25
$ws.Range("J7").Value = 12437.5
$ws.Range("K7").Value = 5642.25
$ws.Range("L7").Value = 12437.5
$ws.Range("M7").Value = -5530.25
$ws.Range("N7").Value = -12661.5
$ws.Range("H40").Value = 5534
$ws.Range("I40").Value = 6732.5
$ws.Range("J40").Value = 4735
$ws.Range("K40").Value = 6732.5
$ws.Range("L40").Value = 4735
$ws.Range("M40").Value = -6596.5
$ws.Range("N40").Value = -5007
$ws.Range("H61").Value = 11899.6
$ws.Range("I61").Value = 8587.666999999999
$ws.Range("J61").Value = 16867.5
$ws.Range("K61").Value = 8587.666999999999
$ws.Range("L61").Value = 16867.5
$ws.Range("M61").Value = -8385.666999999999
$ws.Range("N61").Value = -17271.5
$ws.Range("H68").Value = 1879.0588
$ws.Range("I68").Value = 1929.6
$ws.Range("J68").Value = 1500
$ws.Range("K68").Value = 1929.6
$ws.Range("L68").Value = 1500
$ws.Range("M68").Value = -1180.6
$ws.Range("N68").Value = -2998
$ws.Range("H71").Value = 1879.0588
$ws.Range("I71").Value = 1929.6
$ws.Range("J71").Value = 1500
$ws.Range("K71").Value = 9648
$ws.Range("L71").Value = 7500
$ws.Range("M71").Value = -5904
$ws.Range("N71").Value = -14988
$ws.Range("H87").Value = 45000
$ws.Range("J87").Value = 45000
$ws.Range("L87").Value = 45000
$ws.Range("N87").Value = -47246
$ws.Range("H90").Value = 45000
$ws.Range("J90").Value = 45000
$ws.Range("L90").Value = 135000
$ws.Range("N90").Value = -146232
$ws.Range("H93").Value = 1994.35
$ws.Range("I93").Value = 3058
$ws.Range("J93").Value = 1538.5
$ws.Range("K93").Value = 3058
$ws.Range("L93").Value = 1538.5
$ws.Range("M93").Value = -1810
$ws.Range("N93").Value = -4034.5
$ws.Range("H113").Value = 11899.6
$ws.Range("I113").Value = 8587.666999999999
$ws.Range("J113").Value = 16867.5
$ws.Range("K113").Value = 8587.666999999999
$ws.Range("L113").Value = 16867.5
$ws.Range("M113").Value = -6417.666999999999
$ws.Range("N113").Value = -21207.5
$ws.Range("H122").Value = 4476.237
$ws.Range("I122").Value = 3135.4
$ws.Range("K122").Value = 9406.200000000001
$ws.Range("M122").Value = -6956.200000000001
$ws.Range("H126").Value = 9039.875
$ws.Range("I126").Value = 5642.25
$ws.Range("J126").Value = 12437.5
$ws.Range("K126").Value = 16926.75
$ws.Range("L126").Value = 37312.5
$ws.Range("M126").Value = -14456.75
$ws.Range("N126").Value = -42252.5

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 18398.154
$ws.Range("J62").Value = 27750.4
$ws.Range("L62").Value = 27750.4
$ws.Range("N62").Value = -28998.4
$ws.Range("H65").Value = 18398.154
$ws.Range("J65").Value = 27750.4
$ws.Range("L65").Value = 138752
$ws.Range("N65").Value = -144992
$ws.Range("H96").Value = 2859.32
$ws.Range("I96").Value = 1411.5555
$ws.Range("K96").Value = 1411.5555
$ws.Range("M96").Value = -38.55549999999994
$ws.Range("H107").Value = 523
$ws.Range("I107").Value = 686.25
$ws.Range("J107").Value = 261.8
$ws.Range("K107").Value = 2058.75
$ws.Range("L107").Value = 785.4000000000001
$ws.Range("M107").Value = -138.75
$ws.Range("N107").Value = -4625.4
$ws.Range("H122").Value = 1888.8148
$ws.Range("I122").Value = 1509.4117
$ws.Range("J122").Value = 2533.8
$ws.Range("K122").Value = 4528.2351
$ws.Range("L122").Value = 7601.400000000001
$ws.Range("M122").Value = -2078.2351
$ws.Range("N122").Value = -12501.4
